$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.322.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.681.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.23%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "683.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.679.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.32%  "

$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("E9").Value = "  -5.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.20%  "

$ws.Range("E13").Value = "  -6.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.301.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -10.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.681.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.350.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("E18").Value = "  -1.33%  "

$ws.Range("E19").Value = "  -9.23%  "

$ws.Range("E20").Value = "  -10.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.649"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.827.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.82%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.13%  "

$ws.Range("E30").Value = "  -10.35%  "

$ws.Range("E31").Value = "  -13.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.93%  "

$ws.Range("E33").Value = "  -9.63%  "

$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.648.67"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  -7.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.75%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("E42").Value = "  -10.38%  "

$ws.Range("E43").Value = "  +0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.941"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "165.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.50%  "

$ws.Range("E47").Value = "  -15.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("E49").Value = "  -6.07%  "

$ws.Range("E50").Value = "  -4.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000271"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.59%  "
